$d = $word.ActiveDocument
$dash = [char]0x2013

# ---------------------------------------------------------------------------
# Helper: paste a clone of the "Inputs to the ML method..." paragraph (which
# carries the Arial/20/20 rPr+pPr formatting we need) immediately before the
# paragraph currently at $beforeIndex, then overwrite its text with
# $placeholder so we have an addressable, non-empty range to refine further.
# Returns the 1-based paragraph index of the freshly inserted paragraph
# (always equal to $beforeIndex, since it is inserted right before it).
# ---------------------------------------------------------------------------
function New-ClonedParagraph($beforeIndex, $placeholder) {
    $tmpl = $d.Paragraphs(4).Range
    $tmpl.Copy()
    $target = $d.Paragraphs($beforeIndex).Range
    $start = $target.Start
    $ins = $d.Range($start, $start)
    $ins.Paste()

    $newPara = $d.Paragraphs($beforeIndex).Range
    $len = $newPara.Text.Length
    $noMark = $d.Range($newPara.Start, $newPara.Start + $len - 1)
    $noMark.Text = $placeholder
    return $beforeIndex
}

# ---------------------------------------------------------------------------
# Helper: get the text-only (no paragraph mark) range for paragraph $idx.
# ---------------------------------------------------------------------------
function Get-ParaTextRange($idx) {
    $p = $d.Paragraphs($idx).Range
    $len = $p.Text.Length
    return $d.Range($p.Start, $p.Start + $len - 1)
}

$rPrXml = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1) Blank paragraph right after "Inputs to the ML method are the transformed values"
$idx = New-ClonedParagraph 5 "X"
$r = Get-ParaTextRange $idx
$r.Delete()

# 2) "Base accuracy measure (MSE - ARIMA, ETS, each level)" - highlighted green
$s2 = "Base accuracy measure (MSE " + $dash + " ARIMA, ETS, each level)"
$idx = New-ClonedParagraph 6 $s2
$r = Get-ParaTextRange $idx
$r.Font.HighlightColorIndex = 4

# 3) "MAE calculation"
$idx = New-ClonedParagraph 7 "MAE calculation"

# 4) "Best case 1 and best case 2"
$idx = New-ClonedParagraph 8 "Best case 1 and best case 2"

# 5) "Plot transformed fitted values for high errors ML - box plots across average MAE of time series"
$s5 = "Plot transformed fitted values for high errors ML " + $dash + " box plots across average MAE of time series"
$idx = New-ClonedParagraph 9 $s5

# 6) "Summarise number of observations, no of rolling windows, no of time series and bottom level series, no of transformed inputs"
$s6 = "Summarise number of observations, no of rolling windows, no of time series and bottom level series, no of transformed inputs"
$idx = New-ClonedParagraph 10 $s6

# 7) "Lambda ranges - suitable range"
$s7 = "Lambda ranges " + $dash + " suitable range"
$idx = New-ClonedParagraph 11 $s7

# 8) Blank paragraph
$idx = New-ClonedParagraph 12 "X"
$r = Get-ParaTextRange $idx
$r.Delete()

# 9) "DeepAR, WaveNet - fitted values (clustering and prediction)" with DeepAR/WaveNet
#    marked as spell-check exceptions (proofErr spellStart/spellEnd), matching
#    Word's automatic markup for words it doesn't recognise.
$idx = New-ClonedParagraph 13 "X"
$r = Get-ParaTextRange $idx
$frag = '<w:p ' + $wNs + '>'
$frag += '<w:proofErr w:type="spellStart"/>'
$frag += '<w:r>' + $rPrXml + '<w:t>DeepAR</w:t></w:r>'
$frag += '<w:proofErr w:type="spellEnd"/>'
$frag += '<w:r>' + $rPrXml + '<w:t xml:space="preserve">, </w:t></w:r>'
$frag += '<w:proofErr w:type="spellStart"/>'
$frag += '<w:r>' + $rPrXml + '<w:t>WaveNet</w:t></w:r>'
$frag += '<w:proofErr w:type="spellEnd"/>'
$tail = ' ' + $dash + ' fitted values (clustering and prediction)'
$frag += '<w:r>' + $rPrXml + '<w:t xml:space="preserve">' + $tail + '</w:t></w:r>'
$frag += '</w:p>'
$r.InsertXML($frag)

# ---------------------------------------------------------------------------
# Style change: DefaultParagraphFont becomes semiHidden.
# ---------------------------------------------------------------------------
$style = $d.Styles("Default Paragraph Font")
$style.Hidden = $true

Write-Host "done"
